$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.437.74'
$ws.Range("E2").Value = '  +0.64%  '

$ws.Range("D3").Value = '3.164.84'
$ws.Range("E3").Value = '  -0.69%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '571.85'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.19%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '164.30'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.08%  '

$ws.Range("E7").Value = '  +0.06%  '

$ws.Range("E8").Value = '  -4.10%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.118'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.00%  '

$ws.Range("E10").Value = '  -1.32%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.386'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.57%  '

$ws.Range("D12").Value = '3.715.70'
$ws.Range("E12").Value = '  -0.67%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.128'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.99%  '

$ws.Range("D14").Value = '64.436.48'
$ws.Range("E14").Value = '  +0.45%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '25.35'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.53%  '

$ws.Range("D16").Value = '3.166.33'
$ws.Range("E16").Value = '  -0.55%  '

$ws.Range("E17").Value = '  -2.28%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '408.94'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.63%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.77'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.75%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.28'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.95%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.10'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.12%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.00'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.18%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '68.68'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.88%  '

$ws.Range("E24").Value = '  -1.85%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.485'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.01%  '

$ws.Range("E26").Value = '  -6.48%  '

$ws.Range("E27").Value = '  +1.06%  '

$ws.Range("E28").Value = '  -0.85%  '

$ws.Range("E29").Value = '  -2.33%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '21.25'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.25%  '

$ws.Range("E31").Value = '  -2.16%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.36'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.92%  '

$ws.Range("E33").Value = '  -1.21%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '155.79'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.15%  '

$ws.Range("E35").Value = '  -2.10%  '

$ws.Range("B36").Value = 'Stacks'
$ws.Range("C36").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.70'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.70%  '

$ws.Range("B37").Value = 'Maker'
$ws.Range("C37").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D37").Value = '2.687.96'
$ws.Range("E37").Value = '  -2.28%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '24.07'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.46%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.10'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.52%  '

$ws.Range("E40").Value = '  -3.04%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0621'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.25%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.47'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.72%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0259'
$ws.Range("D43").Style = "Normal"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '21.54'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.68%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '291.60'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.25%  '

$ws.Range("E46").Value = '  +0.03%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0987'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.87%  '

$ws.Range("E48").Value = '  -7.56%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '10.46'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.38%  '

$ws.Range("E50").Value = '  -1.72%  '

$ws.Range("E51").Value = '  -5.72%  '
